$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22
$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "walkingToRunning"
$ws.Range("C22").Value = 38.35754998722376
$ws.Range("D22").Value = -16.75077860103589
$ws.Range("E22").Value = 9.333282041253318
$ws.Range("F22").Value = -9.858623504638672
$ws.Range("G22").Value = -13.49863147735596
$ws.Range("H22").Value = -4.903345584869385

# Row 23
$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "walkingToRunning"
$ws.Range("C23").Value = 4.113532688307743
$ws.Range("D23").Value = -10.04574092130529
$ws.Range("E23").Value = 6.603543145316186
$ws.Range("F23").Value = 1.233975768089294
$ws.Range("G23").Value = 17.44783401489258
$ws.Range("H23").Value = -5.1009521484375

# Row 24
$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "walkingToRunning"
$ws.Range("C24").Value = -28.9259203975992
$ws.Range("D24").Value = 2.342653398928785
$ws.Range("E24").Value = -14.07138639651451
$ws.Range("F24").Value = 6.964565753936768
$ws.Range("G24").Value = -6.484397888183594
$ws.Range("H24").Value = -1.327785611152649

# Row 25
$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "walkingToRunning"
$ws.Range("C25").Value = 42.44476240762218
$ws.Range("D25").Value = -44.82572951820286
$ws.Range("E25").Value = -4.690105651476593
$ws.Range("F25").Value = -2.914164066314697
$ws.Range("G25").Value = -17.45715522766113
$ws.Range("H25").Value = -3.230347871780396

# Row 26
$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "walkingToRunning"
$ws.Range("C26").Value = 47.95408476361651
$ws.Range("D26").Value = -34.91908646222309
$ws.Range("E26").Value = 8.182450048671846
$ws.Range("F26").Value = -0.99935120344162
$ws.Range("G26").Value = 3.104047536849976
$ws.Range("H26").Value = 2.052831172943115

# Row 27
$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "walkingToRunning"
$ws.Range("C27").Value = 16.99495727823363
$ws.Range("D27").Value = -33.59686962269733
$ws.Range("E27").Value = 18.0336912966659
$ws.Range("F27").Value = -10.20110511779785
$ws.Range("G27").Value = -13.58332061767578
$ws.Range("H27").Value = -4.393616676330566

# Row 28
$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "walkingToRunning"
$ws.Range("C28").Value = 0.5081373297654075
$ws.Range("D28").Value = -34.04575231800982
$ws.Range("E28").Value = 16.78914682761511
$ws.Range("F28").Value = 3.081676959991455
$ws.Range("G28").Value = 6.77015495300293
$ws.Range("H28").Value = 9.200222969055176

# Row 29
$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "walkingToRunning"
$ws.Range("C29").Value = -27.83808203513578
$ws.Range("D29").Value = 3.896602891246781
$ws.Range("E29").Value = -12.58596580072954
$ws.Range("F29").Value = 4.065448760986328
$ws.Range("G29").Value = -7.618904590606689
$ws.Range("H29").Value = -0.8031428456306458

# Row 30
$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "walkingToRunning"
$ws.Range("C30").Value = 21.10497084907895
$ws.Range("D30").Value = -35.42225003686769
$ws.Range("E30").Value = -4.233091697929673
$ws.Range("F30").Value = -1.262471556663513
$ws.Range("G30").Value = -17.20468711853027
$ws.Range("H30").Value = -2.038849592208862

# Row 31
$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "walkingToRunning"
$ws.Range("C31").Value = 45.68115871737581
$ws.Range("D31").Value = -33.5533879203098
$ws.Range("E31").Value = 10.18618937901121
$ws.Range("F31").Value = -0.1183775141835212
$ws.Range("G31").Value = 2.852112531661988
$ws.Range("H31").Value = 1.839245676994324

Write-Host "Rows 22-31 added"